$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Remove obsolete rows 58-80 (old entries no longer present)
$ws.Rows("58:80").Delete()

# Step 2: Overwrite rows 48-57 with the updated records
# Row 48 (was case -493)
$ws.Range("A48").Value = "'-493"
$ws.Range("B48").Value = "'6/27/2025"
$ws.Range("C48").Value = "'JUFRE 424"
$ws.Range("D48").Value = "'15"
$ws.Range("E48").Value = "'807817955"
$ws.Range("F48").Value = "'Optical Power"
$ws.Range("G48").Value = "'Pendiente"
$ws.Range("H48").Value = "'Desmontar columna de 168 mm y traspasar redes a columna comunitaria"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = "'Desmonte"
$ws.Range("K48").Value = "'Sin equipos"
$ws.Range("L48").Value = "'Pasante"
$ws.Range("M48").Value = -58.432644
$ws.Range("N48").Value = -34.595434
$ws.Range("O48").Value = "'Palermo"
$ws.Range("P48").Value = "'Capital Sur"

# Row 49 (was case -501)
$ws.Range("A49").Value = "'-501"
$ws.Range("B49").Value = "'7/3/2025"
$ws.Range("C49").Value = "'Cabello 3107"
$ws.Range("D49").Value = "'14"
$ws.Range("E49").Value = "'807971967"
$ws.Range("F49").Value = "'Optical Power"
$ws.Range("G49").Value = "'Pendiente"
$ws.Range("H49").Value = "'Aplomar"
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = "'Aplomo"
$ws.Range("K49").Value = "'Sin equipos"
$ws.Range("L49").Value = "'Terminal"
$ws.Range("M49").Value = -58.405749
$ws.Range("N49").Value = -34.58224
$ws.Range("O49").Value = "'Recoleta"
$ws.Range("P49").Value = "'Capital Sur"

# Row 50 (was case -502)
$ws.Range("A50").Value = "'-502"
$ws.Range("B50").Value = "'7/7/2025"
$ws.Range("C50").Value = "'Tagle 2562"
$ws.Range("D50").Value = "'14"
$ws.Range("E50").Value = "'808036198"
$ws.Range("F50").Value = "'Optical Power"
$ws.Range("G50").Value = "'Pendiente"
$ws.Range("H50").Value = "'Colocar columna para pedir traspaso nodo teco"
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = "'Cambio"
$ws.Range("K50").Value = "'Nodo Teco"
$ws.Range("L50").Value = "'Pasante"
$ws.Range("M50").Value = -58.400188
$ws.Range("N50").Value = -34.583882
$ws.Range("O50").Value = "'Recoleta"
$ws.Range("P50").Value = "'Capital Sur"

# Row 51 (was case -506)
$ws.Range("A51").Value = "'-506"
$ws.Range("B51").Value = "'7/11/2025"
$ws.Range("C51").Value = "'Espinosa 591"
$ws.Range("D51").Value = "'6"
$ws.Range("E51").Value = "'808150511"
$ws.Range("F51").Value = "'Optical Power"
$ws.Range("G51").Value = "'Pendiente"
$ws.Range("H51").Value = "'Picada"
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = "'Cambio"
$ws.Range("K51").Value = "'Nodo Teco"
$ws.Range("L51").Value = "'Pasante"
$ws.Range("M51").Value = -58.449
$ws.Range("N51").Value = -34.616077
$ws.Range("O51").Value = "'Boedo"
$ws.Range("P51").Value = "'Capital Sur"

# Row 52 (was case -511)
$ws.Range("A52").Value = "'-511"
$ws.Range("B52").Value = "'7/14/2025"
$ws.Range("C52").Value = "'Carlos Melo 491"
$ws.Range("D52").Value = "'4"
$ws.Range("E52").Value = "'808194932"
$ws.Range("F52").Value = "'Optical Power"
$ws.Range("G52").Value = "'Pendiente"
$ws.Range("H52").Value = "'Picada"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = "'Cambio"
$ws.Range("K52").Value = "'Sin equipos"
$ws.Range("L52").Value = "'Terminal"
$ws.Range("M52").Value = -58.363292
$ws.Range("N52").Value = -34.642869
$ws.Range("O52").Value = "'San Telmo"
$ws.Range("P52").Value = "'Capital Sur"

# Row 53 (was case -517)
$ws.Range("A53").Value = "'-517"
$ws.Range("B53").Value = "'7/16/2025"
$ws.Range("C53").Value = "'Av Dorrego 2721"
$ws.Range("D53").Value = "'14"
$ws.Range("E53").Value = "'808373635"
$ws.Range("F53").Value = "'Optical Power"
$ws.Range("G53").Value = "'Pendiente"
$ws.Range("H53").Value = "'Cambiar columna 114 base corroida y cable de fo cortado"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = "'Cambio"
$ws.Range("K53").Value = "'Sin equipos"
$ws.Range("L53").Value = "'Pasante"
$ws.Range("M53").Value = -58.432805
$ws.Range("N53").Value = -34.574345
$ws.Range("O53").Value = "'Palermo"
$ws.Range("P53").Value = "'Capital Sur"

# Row 54 (was case -529)
$ws.Range("A54").Value = "'-529"
$ws.Range("B54").Value = "'7/23/2025"
$ws.Range("C54").Value = "'Libertad 820"
$ws.Range("D54").Value = "'1"
$ws.Range("E54").Value = "'ICD30189941"
$ws.Range("F54").Value = "'Optical Power"
$ws.Range("G54").Value = "'Pendiente"
$ws.Range("H54").Value = "'Colocar columna hablar con Pablo si hay dudas"
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = "'Cambio"
$ws.Range("K54").Value = "'Sin equipos"
$ws.Range("L54").Value = "'Pasante"
$ws.Range("M54").Value = -58.384097
$ws.Range("N54").Value = -34.598913
$ws.Range("O54").Value = "'Recoleta"
$ws.Range("P54").Value = "'Capital Sur"

# Row 55 (was case -531)
$ws.Range("A55").Value = "'-531"
$ws.Range("B55").Value = "'7/25/2025"
$ws.Range("C55").Value = "'Joaquin V Gonzalez 4632"
$ws.Range("D55").Value = "'11"
$ws.Range("E55").Value = "'808530239"
$ws.Range("F55").Value = "'Optical Power"
$ws.Range("G55").Value = "'Pendiente"
$ws.Range("H55").Value = "'Cambiar por prfv y usar esa 114 en Libertad 820"
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = "'Cambio"
$ws.Range("K55").Value = "'Sin equipos"
$ws.Range("L55").Value = "'Pasante"
$ws.Range("M55").Value = -58.513643
$ws.Range("N55").Value = -34.594169
$ws.Range("O55").Value = "'Paternal"
$ws.Range("P55").Value = "'Capital Norte"

# Row 56 (was case -593)
$ws.Range("A56").Value = "'-593"
$ws.Range("B56").Value = "'9/10/2025"
$ws.Range("C56").Value = "'Husares 2250"
$ws.Range("D56").Value = "'13"
$ws.Range("E56").Value = "'809642190"
$ws.Range("F56").Value = "'Optical Power"
$ws.Range("G56").Value = "'Pendiente"
$ws.Range("H56").Value = "'Picada"
$ws.Range("I56").Value = 1
$ws.Range("J56").Value = "'Cambio"
$ws.Range("K56").Value = "'Sin equipos"
$ws.Range("L56").Value = "'Pasante"
$ws.Range("M56").Value = -58.443269
$ws.Range("N56").Value = -34.552209
$ws.Range("O56").Value = "'Saavedra"
$ws.Range("P56").Value = "'Capital Norte"

# Row 57 (was case -594)
$ws.Range("A57").Value = "'-594"
$ws.Range("B57").Value = "'9/10/2025"
$ws.Range("C57").Value = "'Vidal 1861"
$ws.Range("D57").Value = "'13"
$ws.Range("E57").Value = "'809642175"
$ws.Range("F57").Value = "'Optical Power"
$ws.Range("G57").Value = "'Pendiente"
$ws.Range("H57").Value = "'Picada"
$ws.Range("I57").Value = 1
$ws.Range("J57").Value = "'Cambio"
$ws.Range("K57").Value = "'Sin equipos"
$ws.Range("L57").Value = "'Pasante"
$ws.Range("M57").Value = -58.458298
$ws.Range("N57").Value = -34.566511
$ws.Range("O57").Value = "'Colegiales"
$ws.Range("P57").Value = "'Capital Norte"

